# Updates the cryptocurrency price ("Price", column D) and hourly
# volume change ("Volume(1h)", column E) columns on Sheet1 to match the
# latest scrape, row by row (rows 2-51).
#
# Column D holds prices formatted as plain text (e.g. "71.980.17" or
# "0.998"), not numbers - some of those strings parse as valid numbers
# (e.g. "0.998", "71.76"), so a naive `.Value = "0.998"` assignment would
# make Excel auto-convert the cell to a genuine number and lose the
# original text formatting/trailing zeros (e.g. "0.550" -> 0.55).
#
# To keep those cells as real text (matching the source workbook, which
# stores them as inline/shared strings) we stage the value in a scratch
# cell (Z1) that is explicitly formatted as Text ("@"), then copy that
# cell and paste-special *values only* into the destination. Pasting
# values-only carries over the already-resolved text value without also
# copying the scratch cell's number format onto the destination cell, so
# the target cells keep their original (default) style.
#
# Values that are NOT valid numbers (e.g. "71.980.17", with two dots, or
# values using the subscript-zero notation "0.0₃0962") are always kept as
# text automatically, so those are just assigned directly.
#
# Column E ("  +0.47%  " etc.) always keeps its leading/trailing spaces
# and "%" sign, so it is never misread as a number and can be assigned
# directly in all cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to stage text-typed numeric-looking values (see above).
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "71.980.17"
$ws.Range("E2").Value = "  +0.47%  "
# Row 3
$ws.Range("D3").Value = "2.697.30"
$ws.Range("E3").Value = "  +2.55%  "
# Row 4
$helper.Value = "0.998"
$helper.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  -0.06%  "
# Row 5
$helper.Value = "596.96"
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -1.50%  "
# Row 6
$helper.Value = "176.10"
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -2.03%  "
# Row 7
$helper.Value = "0.999"
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -0.08%  "
# Row 8
$helper.Value = "0.522"
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -0.69%  "
# Row 9
$ws.Range("D9").Value = "2.695.77"
$ws.Range("E9").Value = "  +2.50%  "
# Row 10
$helper.Value = "0.171"
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +1.96%  "
# Row 11
$helper.Value = "0.169"
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +2.54%  "
# Row 12
$helper.Value = "0.353"
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +1.52%  "
# Row 13
$helper.Value = "5.01"
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -1.09%  "
# Row 14
$ws.Range("D14").Value = "3.180.43"
$ws.Range("E14").Value = "  +1.74%  "
# Row 15
$helper.Value = "0.0000185"
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -0.59%  "
# Row 16
$ws.Range("D16").Value = "71.786.10"
$ws.Range("E16").Value = "  +0.41%  "
# Row 17
$helper.Value = "26.24"
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -1.78%  "
# Row 18
$ws.Range("D18").Value = "2.690.23"
$ws.Range("E18").Value = "  +2.30%  "
# Row 19
$helper.Value = "12.11"
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +5.41%  "
# Row 20
$helper.Value = "8.03"
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +1.70%  "
# Row 21
$helper.Value = "369.69"
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -3.12%  "
# Row 22
$helper.Value = "4.15"
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +0.65%  "
# Row 23
$helper.Value = "2.02"
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +1.23%  "
# Row 24
$helper.Value = "71.76"
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -1.28%  "
# Row 25
$ws.Range("E25").Value = "  -0.07%  "
# Row 26
$helper.Value = "4.33"
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -2.39%  "
# Row 27
$helper.Value = "9.87"
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -1.38%  "
# Row 28
$ws.Range("D28").Value = "2.829.47"
$ws.Range("E28").Value = "  +2.46%  "
# Row 29
$ws.Range("E29").Value = "  -0.06%  "
# Row 30
$ws.Range("D30").Value = "0.0₃0962"
$ws.Range("E30").Value = "  -0.16%  "
# Row 31
$helper.Value = "8.07"
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -0.26%  "
# Row 32
$helper.Value = "507.79"
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -7.86%  "
# Row 33
$ws.Range("E33").Value = "  -3.52%  "
# Row 34
$ws.Range("E34").Value = "  -1.00%  "
# Row 35
$helper.Value = "0.999"
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -0.07%  "
# Row 36
$helper.Value = "162.72"
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -2.12%  "
# Row 37
$helper.Value = "19.43"
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +1.01%  "
# Row 38
$helper.Value = "19.06"
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -0.67%  "
# Row 39
$helper.Value = "1.38"
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -2.31%  "
# Row 40
$helper.Value = "0.109"
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -5.46%  "
# Row 41
$helper.Value = "1.81"
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -3.73%  "
# Row 42
$ws.Range("E42").Value = "  -0.09%  "
# Row 43
$helper.Value = "5.03"
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -0.62%  "
# Row 44
$helper.Value = "2.58"
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -1.75%  "
# Row 45
$helper.Value = "0.333"
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -0.23%  "
# Row 46
$helper.Value = "155.45"
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +1.66%  "
# Row 47
$helper.Value = "39.22"
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -1.04%  "
# Row 48
$helper.Value = "3.72"
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +1.94%  "
# Row 49
$helper.Value = "1.74"
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +3.14%  "
# Row 50
$helper.Value = "0.550"
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +2.53%  "
# Row 51
$helper.Value = "0.0765"
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +0.74%  "

# Clean up the scratch cell so it does not end up in the saved workbook.
$helper.Clear()
